$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.272.61"
$ws.Range("E2").Value = "  -0.45%  "
$ws.Range("D3").Value = "2.389.39"
$ws.Range("E3").Value = "  -3.47%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "549.85"
$ws.Range("E5").Value = "  -0.70%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.07"
$ws.Range("E6").Value = "  -3.33%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.534"
$ws.Range("E8").Value = "  -10.75%  "
$ws.Range("D9").Value = "2.387.16"
$ws.Range("E9").Value = "  -3.53%  "
$ws.Range("E10").Value = "  -2.05%  "
$ws.Range("E11").Value = "  +0.24%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.30"
$ws.Range("E12").Value = "  -3.00%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.53"
$ws.Range("E14").Value = "  -2.50%  "
$ws.Range("D15").Value = "2.819.51"
$ws.Range("E15").Value = "  -3.47%  "
$ws.Range("E16").Value = "  -1.59%  "
$ws.Range("D17").Value = "61.193.69"
$ws.Range("E17").Value = "  -0.32%  "
$ws.Range("D18").Value = "2.387.23"
$ws.Range("E18").Value = "  -3.65%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.76"
$ws.Range("E19").Value = "  -3.89%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.13"
$ws.Range("E20").Value = "  -1.60%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "318.93"
$ws.Range("E21").Value = "  -0.95%  "
$ws.Range("E22").Value = "  -4.76%  "
$ws.Range("E23").Value = "  -0.05%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.92"
$ws.Range("E24").Value = "  +3.31%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "63.62"
$ws.Range("E25").Value = "  -0.64%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.29"
$ws.Range("E26").Value = "  +7.12%  "
$ws.Range("E27").Value = "  +0.15%  "
$ws.Range("D28").Value = "2.505.44"
$ws.Range("E28").Value = "  -3.77%  "
$ws.Range("D29").Value = "0.0₃0928"
$ws.Range("E29").Value = "  -7.06%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "527.28"
$ws.Range("E30").Value = "  -3.86%  "
$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.12"
$ws.Range("E31").Value = "  -2.36%  "
$ws.Range("B32").Value = "Fetch.AI"
$ws.Range("C32").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.43"
$ws.Range("E32").Value = "  -5.15%  "
$ws.Range("E33").Value = "  -3.67%  "
$ws.Range("E34").Value = "  -3.85%  "
$ws.Range("E35").Value = "  -0.23%  "
$ws.Range("E36").Value = "  -0.01%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.53"
$ws.Range("E37").Value = "  -6.47%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.70"
$ws.Range("E38").Value = "  -3.69%  "
$ws.Range("E39").Value = "  -1.96%  "
$ws.Range("E40").Value = "  +5.78%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "18.10"
$ws.Range("E41").Value = "  -2.39%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "139.50"
$ws.Range("E42").Value = "  -3.89%  "
$ws.Range("E43").Value = "  +0.07%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "40.29"
$ws.Range("E44").Value = "  -0.31%  "
$ws.Range("E45").Value = "  -9.40%  "
$ws.Range("E46").Value = "  -0.45%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "140.85"
$ws.Range("E47").Value = "  -4.60%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "20.11"
$ws.Range("E48").Value = "  -6.08%  "
$ws.Range("E49").Value = "  -3.82%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.576"
$ws.Range("E50").Value = "  -3.38%  "
$ws.Range("E51").Value = "  -1.14%  "
